$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 4,23
$arr[0,0] = 0.979323308270677
$arr[0,1] = 0.0075187969924812
$arr[0,2] = 0.0093984962406015
$arr[0,3] = 0.0507518796992481
$arr[0,4] = 0.0056390977443609
$arr[0,5] = 0.0093984962406015
$arr[0,6] = 0.973684210526316
$arr[0,7] = 0.0244360902255639
$arr[0,8] = 0.951127819548872
$arr[0,9] = 0.964285714285714
$arr[0,10] = 0
$arr[0,11] = 0.0037593984962406
$arr[0,12] = 0.99812030075188
$arr[0,13] = 0
$arr[0,14] = 0.996240601503759
$arr[0,15] = 0
$arr[0,16] = 0.890977443609023
$arr[0,17] = 0.0526315789473684
$arr[0,18] = 0.0639097744360902
$arr[0,19] = 0.994360902255639
$arr[0,20] = 0.0056390977443609
$arr[0,21] = 0.0018796992481203
$arr[0,22] = 0.018796992481203
$arr[1,0] = 0.0056390977443609
$arr[1,1] = 0.0244360902255639
$arr[1,2] = 0.0037593984962406
$arr[1,3] = 0.0056390977443609
$arr[1,4] = 0.0037593984962406
$arr[1,5] = 0.977443609022556
$arr[1,6] = 0.0093984962406015
$arr[1,7] = 0.0131578947368421
$arr[1,8] = 0
$arr[1,9] = 0.0018796992481203
$arr[1,10] = 0.93609022556391
$arr[1,11] = 0.0056390977443609
$arr[1,12] = 0
$arr[1,13] = 0.99812030075188
$arr[1,14] = 0.0037593984962406
$arr[1,15] = 1
$arr[1,16] = 0.093984962406015
$arr[1,17] = 0.934210526315789
$arr[1,18] = 0.926691729323308
$arr[1,19] = 0
$arr[1,20] = 0.0037593984962406
$arr[1,21] = 0
$arr[1,22] = 0.0075187969924812
$arr[2,0] = 0.0037593984962406
$arr[2,1] = 0.0075187969924812
$arr[2,2] = 0
$arr[2,3] = 0.941729323308271
$arr[2,4] = 0.988721804511278
$arr[2,5] = 0.0112781954887218
$arr[2,6] = 0.0056390977443609
$arr[2,7] = 0.949248120300752
$arr[2,8] = 0.0093984962406015
$arr[2,9] = 0.0319548872180451
$arr[2,10] = 0.0018796992481203
$arr[2,11] = 0
$arr[2,12] = 0.0018796992481203
$arr[2,13] = 0
$arr[2,14] = 0
$arr[2,15] = 0
$arr[2,16] = 0.0112781954887218
$arr[2,17] = 0.0037593984962406
$arr[2,18] = 0.0018796992481203
$arr[2,19] = 0.0056390977443609
$arr[2,20] = 0.984962406015038
$arr[2,21] = 0.996240601503759
$arr[2,22] = 0.969924812030075
$arr[3,0] = 0.0112781954887218
$arr[3,1] = 0.960526315789474
$arr[3,2] = 0.986842105263158
$arr[3,3] = 0.0018796992481203
$arr[3,4] = 0.0018796992481203
$arr[3,5] = 0.0018796992481203
$arr[3,6] = 0.0112781954887218
$arr[3,7] = 0.0131578947368421
$arr[3,8] = 0.037593984962406
$arr[3,9] = 0
$arr[3,10] = 0.0601503759398496
$arr[3,11] = 0.990601503759398
$arr[3,12] = 0
$arr[3,13] = 0.0018796992481203
$arr[3,14] = 0
$arr[3,15] = 0
$arr[3,16] = 0.0037593984962406
$arr[3,17] = 0.0075187969924812
$arr[3,18] = 0.0075187969924812
$arr[3,19] = 0
$arr[3,20] = 0.0056390977443609
$arr[3,21] = 0.0018796992481203
$arr[3,22] = 0.0018796992481203

$ws.Range("B2:X5").Value = $arr
